# Auto commit at 2025-10-02  8:20:38.51
#
# Extends the daily-stats sheet with two new rows (62 and 63) for the
# 2025-10-01 reading, mirroring the pattern used by the previous rows,
# and converts the D56:D61 "minutes -> Excel time" formulas into one
# shared-formula group (as Excel does when a formula is filled across
# a contiguous range), then moves the active selection to J63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-enter the D56:D61 formula as a single fill across the range so ---
# --- the engine records it as one shared formula group (si="0").      ---
$ws.Range("D56:D61").Formula = "=C56/(24*60)"

# --- Row 62: 2025-10-01, 四方坪站 ---
$ws.Range("A62").Value = 45931
$ws.Range("B62").Value = "四方坪站"
$ws.Range("C62").Formula = "=16799/127"
$ws.Range("D62").Formula = "=C62/(24*60)"
$ws.Range("E62").Formula = "=10368.41/127"
$ws.Range("F62").Formula = "=3574.11/127"
$ws.Range("G62").Formula = "=10368.41/(16799/60)"
$ws.Range("H62").Formula = "=436/127"

# --- Row 63: 2025-10-01, 高岭站 ---
$ws.Range("A63").Value = 45931
$ws.Range("B63").Value = "高岭站"
$ws.Range("C63").Formula = "=6591/36"
$ws.Range("D63").Formula = "=C63/(24*60)"
$ws.Range("E63").Formula = "=4791.76/36"
$ws.Range("F63").Formula = "=1184.73/36"
$ws.Range("G63").Formula = "=4791.76/(6591/60)"
$ws.Range("H63").Formula = "=178/36"

# --- Move the active selection from J61 to J63 ---
$ws.Range("J63").Select()
